# "#5: property boat&car done"
# Fill in the real data row for the "汽車" (car) property sheet, which
# previously only held a stray placeholder row (A1:G2), expanding it to the
# full normalized record layout (A1:N2) used by the other property sheets
# (land/building/stock/...): name, capacity, owner, register_date,
# register_reason, acquire_value, property_category, category, date,
# legislator_name, legislator_id, source_file, index.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)   # "汽車" sheet (3rd tab)

# ---- Header row (row 1) ----------------------------------------------
$ws3.Range("B1").Value = "name"
$ws3.Range("C1").Value = "capacity"
$ws3.Range("D1").Value = "owner"
$ws3.Range("E1").Value = "register_date"
$ws3.Range("F1").Value = "register_reason"
$ws3.Range("G1").Value = "acquire_value"
$ws3.Range("H1").Value = "property_category"
$ws3.Range("I1").Value = "category"
$ws3.Range("J1").Value = "date"
$ws3.Range("K1").Value = "legislator_name"
$ws3.Range("L1").Value = "legislator_id"
$ws3.Range("M1").Value = "source_file"
$ws3.Range("N1").Value = "index"

# Apply the same (bold / bordered) header formatting used by B1:G1 to the
# newly added H1:N1 header cells, without disturbing their values.
$ws3.Range("B1:G1").Copy() | Out-Null
$ws3.Range("H1:N1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# ---- Data row (row 2) --------------------------------------------------
$ws3.Range("A2").Value = 34
$ws3.Range("B2").Value = "曰產"
$ws3.Range("C2").Value = 1998
$ws3.Range("D2").Value = "周桂香"
$ws3.Range("E2").Value = "92年10月15日"
$ws3.Range("F2").Value = "買賣"
$ws3.Range("G2").Value = "(超過五年）"
$ws3.Range("H2").Value = "land"
$ws3.Range("I2").Value = "normal"
$ws3.Range("J2").Value = "2013-12-11"
$ws3.Range("K2").Value = "許智傑"
$ws3.Range("L2").Value = 1750
$ws3.Range("M2").Value = "tmpd3cb1"
$ws3.Range("N2").Value = 34

# Apply the same plain data formatting used by B2:G2 to the newly added
# H2:N2 data cells, without disturbing their values.
$ws3.Range("B2:G2").Copy() | Out-Null
$ws3.Range("H2:N2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
